$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = -0.999490387077816
$ws.Range("C2").Value = -167965926.03901
$ws.Range("D2").Value = 0.926355130140787
$ws.Range("E2").Value = 0.999489979962112
$ws.Range("F2").Value = -0.926364192749595
$ws.Range("G2").Value = 100753.953522504
$ws.Range("H2").Value = 7053899.19618733
$ws.Range("I2").Value = -2399153.90096887
$ws.Range("J2").Value = 0.929257128540929
$ws.Range("K2").Value = 0.999622989973222
$ws.Range("L2").Value = -0.926364192749595
$ws.Range("M2").Value = 41.9784639654396
$ws.Range("N2").Value = 2930.17188190043
$ws.Range("O2").Value = -2399153.90096887

# Row 3
$ws.Range("B3").Value = 0.958630856791998
$ws.Range("C3").Value = 208.691221224262
$ws.Range("D3").Value = -0.982270111965695
$ws.Range("E3").Value = -0.960485486419506
$ws.Range("F3").Value = 0.977776237804515
$ws.Range("G3").Value = -11.6486991963556
$ws.Range("H3").Value = -739.100551819794
$ws.Range("I3").Value = 3.28039420639715
$ws.Range("J3").Value = -0.981147411280848
$ws.Range("K3").Value = -0.858674536582941
$ws.Range("L3").Value = 0.977776237804515
$ws.Range("M3").Value = -0.0494685352528135
$ws.Range("N3").Value = -2.80924484210361
$ws.Range("O3").Value = 3.28039420639715

# Row 5
$ws.Range("B5").Value = 0.959724684595245
$ws.Range("C5").Value = 210.937762444454
$ws.Range("D5").Value = -0.981720824377626
$ws.Range("E5").Value = -0.962931925400802
$ws.Range("F5").Value = 0.977548068821597
$ws.Range("G5").Value = -11.9543803291666
$ws.Range("H5").Value = -760.853197313738
$ws.Range("I5").Value = 3.31115547025337
$ws.Range("J5").Value = -0.997005467917682
$ws.Range("K5").Value = -0.945219214524355
$ws.Range("L5").Value = 0.977548068821597
$ws.Range("M5").Value = -0.0737430148437234
$ws.Range("N5").Value = -4.53651247098088
$ws.Range("O5").Value = 3.31115547025337
